# Cronograma de proyecto.xlsx — update monthly tracking marks (rows 16-24)
# Several cells previously marked with the text "x" are updated: some move
# to a numeric 1 in a different month column, others are simply cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# Row 16: mark moves from column N ("x") to column M (numeric 1)
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = ""

# Row 17: "x" in N becomes numeric 1 in N
$ws.Range("N17").Value = 1

# Row 18: "x" in N becomes numeric 1 in N
$ws.Range("N18").Value = 1

# Row 19: "x" in N becomes numeric 1 in N
$ws.Range("N19").Value = 1

# Row 20: "x" in N becomes numeric 1 in N
$ws.Range("N20").Value = 1

# Row 21: mark moves from column N ("x") to column O (numeric 1)
$ws.Range("N21").Value = ""
$ws.Range("O21").Value = 1

# Row 22: "x" in N is cleared (O22 keeps its own "x")
$ws.Range("N22").Value = ""

# Row 23: "x" in N is cleared (O23 keeps its own "x")
$ws.Range("N23").Value = ""

# Row 24: "x" in N is cleared (O24 keeps its own "x")
$ws.Range("N24").Value = ""
